# Automatische test-sync: 2025-07-27 19:23:50
#
# Adds the new "Productinformatie" test-mail row to the Logs sheet,
# rolls the matching category tally into the Dashboard sheet, and
# extends the conditional formatting ranges + the dashboard bar chart
# so everything keeps covering the newly-added row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append row 8 with the new test mail entry.
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A8").Value = "Hebben we EcoPro-700 nog op voorraad?"
$logs.Range("B8").Value = "mailmind.test@zohomail.eu"
$logs.Range("C8").Value = "Testmail #6: Hebben we EcoPro-700 nog op voorraad?"
$logs.Range("D8").Value = "Productinformatie"
$logs.Range("E8").Value = "Beste klant,`nBedankt voor uw e-mail. Op dit moment hebben we nog voorraad van de EcoPro-700. U kunt deze bestellen via onze website of contact opnemen met onze verkoopafdeling voor meer informatie.`nMet vriendelijke groet,`n[Bedrijfsnaam] - Verkoopafdeling"
$logs.Range("F8").Value = "2025-07-27 19:23:23"
$logs.Range("G8").Value = "Ja"
$logs.Range("H8").Value = "Nee"
$logs.Range("I8").Value = "Ja"
$logs.Range("J8").Value = "Nee"

# ---------------------------------------------------------------
# 2. Extend the conditional formatting sqref ranges from row 7 to
#    row 8 for every column that carries a rule.
# ---------------------------------------------------------------
$logs.Range("D2:D7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D8"))
$logs.Range("G2:G7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G8"))
$logs.Range("H2:H7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H8"))
$logs.Range("I2:I7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I8"))
$logs.Range("J2:J7").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J8"))

# ---------------------------------------------------------------
# 3. Dashboard sheet: add the "Productinformatie" tally row.
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A5").Value = "Productinformatie"
$dash.Range("B5").Value = 1

# ---------------------------------------------------------------
# 4. Update the bar chart so its category/value series cover the
#    newly added dashboard row (A2:A5 / B2:B5 instead of A2:A4 / B2:B4).
# ---------------------------------------------------------------
$chart = $dash.ChartObjects().Item(1).Chart
$series = $chart.SeriesCollection().Item(1)
$series.Formula = "=SERIES(Dashboard!`$B`$1,Dashboard!`$A`$2:`$A`$5,Dashboard!`$B`$2:`$B`$5,1)"
